# Auto-generated edit script for cryptos.xlsx price/volume refresh
# (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to stay a text string (matches the source sheet,
    # which stores these as inline/shared strings, not numbers),
    # then restore the default "Normal" style so no stray number
    # format sticks around on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '63.894.82'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '3.498.26'
$ws.Range('E3').Value = '  -2.47%  '
Set-TextCell $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  +0.02%  '
Set-TextCell $ws.Range('D5') '584.38'
$ws.Range('E5').Value = '  -3.07%  '
Set-TextCell $ws.Range('D6') '131.83'
$ws.Range('E6').Value = '  -4.31%  '
$ws.Range('D7').Value = '3.500.60'
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('E8').Value = '  +0.07%  '
Set-TextCell $ws.Range('D9') '0.489'
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -1.32%  '
Set-TextCell $ws.Range('D11') '7.09'
$ws.Range('E11').Value = '  -1.71%  '
Set-TextCell $ws.Range('D12') '0.386'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').Value = '4.069.70'
$ws.Range('E13').Value = '  -3.07%  '
Set-TextCell $ws.Range('D14') '27.84'
$ws.Range('E14').Value = '  -1.15%  '
Set-TextCell $ws.Range('D15') '0.0000179'
$ws.Range('E15').Value = '  -4.27%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '3.478.34'
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').Value = '63.982.76'
$ws.Range('E18').Value = '  -2.72%  '
Set-TextCell $ws.Range('D19') '10.01'
$ws.Range('E19').Value = '  -0.01%  '
Set-TextCell $ws.Range('D20') '14.34'
$ws.Range('E20').Value = '  -1.89%  '
Set-TextCell $ws.Range('D21') '5.66'
$ws.Range('E21').Value = '  -3.94%  '
Set-TextCell $ws.Range('D22') '391.47'
$ws.Range('E22').Value = '  -1.43%  '
Set-TextCell $ws.Range('D23') '0.578'
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('D24').Value = '3.625.18'
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws.Range('D25') '72.80'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws.Range('D26') '0.999'
$ws.Range('E26').Value = '  -0.02%  '
Set-TextCell $ws.Range('D27') '0.0000111'
$ws.Range('E27').Value = '  -6.28%  '
Set-TextCell $ws.Range('D28') '1.57'
$ws.Range('E28').Value = '  -4.18%  '
Set-TextCell $ws.Range('D29') '0.999'
$ws.Range('E29').Value = '  -0.08%  '
Set-TextCell $ws.Range('D30') '7.43'
$ws.Range('E30').Value = '  -8.67%  '
$ws.Range('E31').Value = '  -6.13%  '
Set-TextCell $ws.Range('D32') '8.21'
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('D33').Value = '3.490.89'
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('E34').Value = '  +0.01%  '
Set-TextCell $ws.Range('D35') '23.78'
$ws.Range('E35').Value = '  -2.82%  '
Set-TextCell $ws.Range('D36') '0.144'
$ws.Range('E36').Value = '  -2.76%  '
Set-TextCell $ws.Range('D37') '5.35'
$ws.Range('E37').Value = '  -1.20%  '
Set-TextCell $ws.Range('D38') '6.95'
$ws.Range('E38').Value = '  -1.73%  '
Set-TextCell $ws.Range('D39') '1.57'
$ws.Range('E39').Value = '  -3.30%  '
Set-TextCell $ws.Range('D40') '168.69'
$ws.Range('E40').Value = '  -0.36%  '
Set-TextCell $ws.Range('D41') '0.0808'
$ws.Range('E41').Value = '  -3.53%  '
Set-TextCell $ws.Range('D42') '0.811'
$ws.Range('E42').Value = '  -3.64%  '
Set-TextCell $ws.Range('D43') '26.03'
$ws.Range('E43').Value = '  -4.02%  '
Set-TextCell $ws.Range('D44') '0.997'
$ws.Range('E44').Value = '  -0.19%  '
Set-TextCell $ws.Range('D45') '41.78'
$ws.Range('E45').Value = '  -3.10%  '
Set-TextCell $ws.Range('D46') '1.20'
$ws.Range('E46').Value = '  -6.27%  '
Set-TextCell $ws.Range('D47') '4.36'
$ws.Range('E47').Value = '  -4.14%  '
Set-TextCell $ws.Range('D48') '1.64'
$ws.Range('E48').Value = '  -3.93%  '
Set-TextCell $ws.Range('D49') '6.88'
$ws.Range('E49').Value = '  -2.21%  '
$ws.Range('D50').Value = '2.437.72'
$ws.Range('E50').Value = '  -0.70%  '
Set-TextCell $ws.Range('D51') '0.0268'
$ws.Range('E51').Value = '  -1.17%  '
